$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 4.85
$ws.Range("J2").Value = 2.7
$ws.Range("K2").Value = 1.78
$ws.Range("L2").Value = 5.7
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 4.4
$ws.Range("O2").Value = 1.72
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.32
$ws.Range("S2").Value = 5.8
$ws.Range("V2").Value = 2.02
$ws.Range("W2").Value = 2.55
$ws.Range("X2").Value = 1.45
$ws.Range("Y2").Value = 4.45
$ws.Range("Z2").Value = 7.4
$ws.Range("AA2").Value = 10
$ws.Range("AB2").Value = 17.5
$ws.Range("AC2").Value = 23
$ws.Range("AE2").Value = 4.4
$ws.Range("AF2").Value = 6
$ws.Range("AG2").Value = 26
$ws.Range("AH2").Value = 250
$ws.Range("AI2").Value = 8
$ws.Range("AJ2").Value = 25
$ws.Range("AK2").Value = 19
$ws.Range("AL2").Value = 110
$ws.Range("AM2").Value = 90
$ws.Range("AN2").Value = 120

# Row 3 updates
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.75
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 2.5
$ws.Range("L3").Value = 3.6
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 3.6
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 10.5
$ws.Range("AA3").Value = 8.5
$ws.Range("AB3").Value = 17.5
$ws.Range("AC3").Value = 14.5
$ws.Range("AD3").Value = 22
$ws.Range("AF3").Value = 7.4
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 18.5
$ws.Range("AK3").Value = 11.5
$ws.Range("AL3").Value = 40
$ws.Range("AM3").Value = 25
$ws.Range("AN3").Value = 29

# Row 7 updates
$ws.Range("G7").Value = 1.91
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 3.7
$ws.Range("J7").Value = 2.47
$ws.Range("K7").Value = 2.18
$ws.Range("L7").Value = 4.1
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 3.35
$ws.Range("T7").Value = 1.28
$ws.Range("U7").Value = 1.38
$ws.Range("V7").Value = 2.8
$ws.Range("W7").Value = 1.87
$ws.Range("X7").Value = 1.83
$ws.Range("Y7").Value = 6.7
$ws.Range("Z7").Value = 8.5
$ws.Range("AA7").Value = 8.5
$ws.Range("AB7").Value = 16
$ws.Range("AC7").Value = 16
$ws.Range("AE7").Value = 7
$ws.Range("AF7").Value = 6.7
$ws.Range("AG7").Value = 16
$ws.Range("AH7").Value = 80
$ws.Range("AI7").Value = 10
$ws.Range("AJ7").Value = 19
$ws.Range("AK7").Value = 12.5
$ws.Range("AL7").Value = 50
$ws.Range("AM7").Value = 35
$ws.Range("AN7").Value = 45

$wb.Save()
